$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column A (category strings) so the shared-strings table entries for them get dropped,
# letting the new header strings claim the lower indices, matching the original export order.
$ws.Range("A2:A7").ClearContents()

# Add new header cells L1:N1, copying header style from K1
$ws.Range("K1").Copy($ws.Range("L1:N1"))
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Re-populate column A (category strings), now registered after the new headers
$ws.Range("A2").Value = 'aon'
$ws.Range("A3").Value = 'aon'
$ws.Range("A4").Value = 'flex'
$ws.Range("A5").Value = 'flex'
$ws.Range("A6").Value = 'sub'
$ws.Range("A7").Value = 'sub'

# Update E/F (percent-like) columns: values now expressed *100
$ws.Range("E2").Value = 62.84644194756554
$ws.Range("F2").Value = 59.95232419547079
$ws.Range("E3").Value = 37.15355805243446
$ws.Range("F3").Value = 65.92741935483872
$ws.Range("E4").Value = 53.40599455040872
$ws.Range("F4").Value = 90.56122448979592
$ws.Range("E5").Value = 46.59400544959128
$ws.Range("F5").Value = 98.39181286549707
$ws.Range("E6").Value = 89.61988304093568
$ws.Range("F6").Value = 20.55464926590538
$ws.Range("E7").Value = 10.38011695906433
$ws.Range("F7").Value = 36.61971830985916

# Add new data columns L, M, N
$ws.Range("L2").Value = 94.10714337626324
$ws.Range("M2").Value = 154974
$ws.Range("N2").Value = 308.0994035785288
$ws.Range("L3").Value = 87.3015904523896
$ws.Range("M3").Value = 108579
$ws.Range("N3").Value = 332.045871559633
$ws.Range("L4").Value = 86.43287869528324
$ws.Range("M4").Value = 88772
$ws.Range("N4").Value = 125.030985915493
$ws.Range("L5").Value = 93.05249603932496
$ws.Range("M5").Value = 114874
$ws.Range("N5").Value = 170.6894502228826
$ws.Range("L6").Value = 18.6886320009998
$ws.Range("M6").Value = 1789
$ws.Range("N6").Value = 14.1984126984127
$ws.Range("L7").Value = 23.27683795941807
$ws.Range("M7").Value = 419
$ws.Range("N7").Value = 16.11538461538462
